# Updated cryptos list - apply per-row price / volume(1h) changes,
# plus the Monero/Cosmos and FTXToken/ARBITRUM row swaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $B, $C, $D, $E) {
    if ($B -ne $null) { $ws.Cells.Item($Row, 2).Value = $B }
    if ($C -ne $null) { $ws.Cells.Item($Row, 3).Value = $C }
    if ($D -ne $null) {
        # Force text format so numeric-looking strings (e.g. "227.30", "0.613")
        # keep their exact original formatting instead of being parsed as numbers.
        $ws.Cells.Item($Row, 4).NumberFormat = "@"
        $ws.Cells.Item($Row, 4).Value = $D
    }
    if ($E -ne $null) { $ws.Cells.Item($Row, 5).Value = $E }
}

# Row 2 - Bitcoin
Set-Row 2 $null $null "38.750.18" "  +0.34%  "
# Row 3 - Ethereum
Set-Row 3 $null $null "2.100.67" "  +0.39%  "
# Row 4 - TetherUSD
Set-Row 4 $null $null $null "  -0.06%  "
# Row 5 - BNB
Set-Row 5 $null $null "227.30" "  -0.62%  "
# Row 6 - XRP
Set-Row 6 $null $null "0.613" "  -0.12%  "
# Row 7 - Solana
Set-Row 7 $null $null "62.38" "  +1.97%  "
# Row 8 - USDC
Set-Row 8 $null $null $null "  -0.05%  "
# Row 9 - Cardano
Set-Row 9 $null $null $null "  +1.89%  "
# Row 10 - Dogecoin
Set-Row 10 $null $null "0.0841" "  -0.14%  "
# Row 11 - TRON
Set-Row 11 $null $null $null "  -1.31%  "
# Row 12 - Chainlink
Set-Row 12 $null $null $null "  +5.64%  "
# Row 13 - WrappedliquidstakedEther2.0
Set-Row 13 $null $null "2.412.55" "  +0.59%  "
# Row 14 - Avalanche
Set-Row 14 $null $null "22.02" "  -1.34%  "
# Row 15 - Polygon
Set-Row 15 $null $null $null "  +3.70%  "
# Row 16 - Polkadot
Set-Row 16 $null $null $null "  +1.12%  "
# Row 17 - WrappedEther
Set-Row 17 $null $null "2.199.35" "  +4.87%  "
# Row 18 - WrappedBTC
Set-Row 18 $null $null "38.738.38" "  +0.58%  "
# Row 19 - Uniswap
Set-Row 19 $null $null $null "  +0.81%  "
# Row 20 - Litecoin
Set-Row 20 $null $null "71.64" "  +0.92%  "
# Row 21 - ShibaInu
Set-Row 21 $null $null $null "  +0.58%  "
# Row 22 - BitcoinCash
Set-Row 22 $null $null "227.95" "  +0.74%  "
# Row 23 - Dai
Set-Row 23 $null $null $null "  +0.00%  "
# Row 24 - Toncoin
Set-Row 24 $null $null "2.36" "  -2.94%  "
# Row 25 - PancakeSwap
Set-Row 25 $null $null "2.30" "  -0.46%  "

# Rows 26 & 27 - Monero and Cosmos swap places (Cosmos now ranks above Monero)
Set-Row 26 "Cosmos" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom" "9.64" "  +1.99%  "
Set-Row 27 "Monero" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" "171.77" "  +0.86%  "

# Row 28 - Kaspa
Set-Row 28 $null $null $null "  +0.77%  "
# Row 29 - ImmutableX
Set-Row 29 $null $null $null "  +3.98%  "
# Row 30 - EthereumClassic
Set-Row 30 $null $null "19.31" "  +0.99%  "
# Row 31 - WEMIXToken
Set-Row 31 $null $null $null "  +7.68%  "
# Row 32 - Stellar
Set-Row 32 $null $null $null "  +0.12%  "
# Row 33 - Filecoin
Set-Row 33 $null $null $null "  +1.11%  "
# Row 34 - InternetComputer(DFINITY)
Set-Row 34 $null $null "4.73" "  -0.30%  "
# Row 35 - THORChain
Set-Row 35 $null $null $null "  +7.81%  "
# Row 36 - Hedera
Set-Row 36 $null $null "0.0616" "  +1.60%  "
# Row 37 - LidoDAOToken
Set-Row 37 $null $null $null "  +0.23%  "
# Row 38 - RenderToken
Set-Row 38 $null $null "3.53" "  -0.85%  "
# Row 39 - BinanceUSD
Set-Row 39 $null $null "1.00" "  -0.02%  "
# Row 40 - InjectiveProtocol
Set-Row 40 $null $null "18.06" "  -2.33%  "
# Row 41 - Aave
Set-Row 41 $null $null "102.70" "  +2.80%  "
# Row 42 - VeChain
Set-Row 42 $null $null $null "  +3.01%  "
# Row 43 - Maker
Set-Row 43 $null $null "1.527.85" "  -1.15%  "
# Row 44 - TrustWalletToken
Set-Row 44 $null $null $null "  +6.55%  "
# Row 45 - FraxShare
Set-Row 45 $null $null "7.88" "  +0.96%  "
# Row 46 - HuobiToken
Set-Row 46 $null $null $null "  -0.79%  "
# Row 47 - Cronos
Set-Row 47 $null $null "0.0914" "  -1.24%  "

# Rows 48 & 49 - FTXToken and ARBITRUM swap places (ARBITRUM now ranks above FTXToken)
Set-Row 48 "ARBITRUM" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb" "1.07" "  +3.13%  "
Set-Row 49 "FTXToken" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt" "4.14" "  -0.73%  "

# Row 50 - MXToken
Set-Row 50 $null $null $null "  -0.55%  "
# Row 51 - RocketPoolETH
Set-Row 51 $null $null "2.298.89" "  +0.41%  "
